$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet (and its _FilterDatabase defined name reference
#    updates automatically with it).
$ws.Name = "2024-04-07"

# 2. Stamp an explicit 15.75pt custom row height on every row from 1
#    through 1000 (matches the re-export done by Google Sheets: rows
#    1-162 and 164-165 already had data, 163 and 166-1000 are blank
#    trailing rows that still get the explicit height written out).
for ($r = 1; $r -le 1000; $r++) {
    $ws.Rows.Item($r).RowHeight = 15.75
}

Write-Output "done"
